$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Range("H18").Value = 1835.8
$ws.Range("I18").Value = 1835.8
$ws.Range("M18").Value = -1551.8
$ws.Range("K18").Value = 1835.8
# Row 28
$ws.Range("I28").Value = 6300.5
$ws.Range("H28").Value = 17716.834
$ws.Range("J28").Value = 29133.166
$ws.Range("N28").Value = -30103.166
$ws.Range("K28").Value = 6300.5
$ws.Range("M28").Value = -5815.5
$ws.Range("L28").Value = 29133.166
# Row 113
$ws.Range("L113").Value = 7164
$ws.Range("N113").Value = -13672
$ws.Range("H113").Value = 7318.467
$ws.Range("J113").Value = 7164
# Row 137
$ws.Range("L137").Value = 11703.6
$ws.Range("H137").Value = 3141.0454
$ws.Range("J137").Value = 3901.2
$ws.Range("N137").Value = -16803.6

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("I2").Value = 847.78723
$ws.Range("H2").Value = 901.86
$ws.Range("M2").Value = -734.78723
$ws.Range("K2").Value = 847.78723
# Row 5
$ws.Range("L5").Value = 0
$ws.Range("I5").Value = 700.5
$ws.Range("H5").Value = 700.5
$ws.Range("J5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("M5").Value = -588.5
$ws.Range("K5").Value = 700.5
# Row 32
$ws.Range("K32").Value = 2422.9714
$ws.Range("I32").Value = 2422.9714
$ws.Range("H32").Value = 4473.1025
$ws.Range("M32").Value = -2135.9714
# Row 92
$ws.Range("H92").Value = 15275
$ws.Range("J92").Value = 15275
$ws.Range("N92").Value = -20267
$ws.Range("L92").Value = 15275
# Row 116
$ws.Range("I116").Value = 847.78723
$ws.Range("H116").Value = 901.86
$ws.Range("K116").Value = 847.78723
$ws.Range("M116").Value = 1446.21277
# Row 122
$ws.Range("N122").Value = -11196.1432
$ws.Range("H122").Value = 2032.2222
$ws.Range("M122").Value = -2948.5
$ws.Range("J122").Value = 2098.7144
$ws.Range("K122").Value = 5398.5
$ws.Range("L122").Value = 6296.1432
$ws.Range("I122").Value = 1799.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("M3").Value = -733.78723
$ws.Range("H3").Value = 901.86
$ws.Range("K3").Value = 847.78723
$ws.Range("I3").Value = 847.78723
# Row 4
$ws.Range("K4").Value = 700.5
$ws.Range("L4").Value = 0
$ws.Range("I4").Value = 700.5
$ws.Range("H4").Value = 700.5
$ws.Range("N4").ClearContents()
$ws.Range("M4").Value = -585.5
$ws.Range("J4").Value = 0
# Row 82
$ws.Range("L82").Value = 39837.5
$ws.Range("H82").Value = 18547.223
$ws.Range("J82").Value = 39837.5
$ws.Range("N82").Value = -40603.5
# Row 85
$ws.Range("L85").Value = 39837.5
$ws.Range("H85").Value = 18547.223
$ws.Range("J85").Value = 39837.5
$ws.Range("N85").Value = -42489.5
# Row 107
$ws.Range("H107").Value = 2360
$ws.Range("J107").Value = 3425
$ws.Range("N107").Value = -7265
$ws.Range("K107").Value = 1827.5
$ws.Range("M107").Value = 92.5
$ws.Range("L107").Value = 3425
$ws.Range("I107").Value = 1827.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 10
$ws.Range("L10").Value = 1625
$ws.Range("I10").Value = 62.666668
$ws.Range("H10").Value = 1234.4166
$ws.Range("J10").Value = 1625
$ws.Range("N10").Value = -1903
$ws.Range("M10").Value = 76.333332
$ws.Range("K10").Value = 62.666668
# Row 17
$ws.Range("K17").Value = 0
$ws.Range("H17").Value = 999
$ws.Range("I17").Value = 0
$ws.Range("M17").ClearContents()
# Row 31
$ws.Range("I31").Value = 2583.1428
$ws.Range("H31").Value = 5013.7
$ws.Range("K31").Value = 2583.1428
$ws.Range("M31").Value = -2288.1428
# Row 34
$ws.Range("I34").Value = 2583.1428
$ws.Range("H34").Value = 5013.7
$ws.Range("M34").Value = -2381.1428
$ws.Range("K34").Value = 2583.1428
# Row 41
$ws.Range("L41").Value = 110000
$ws.Range("H41").Value = 110000
$ws.Range("J41").Value = 110000
$ws.Range("N41").Value = -110856
# Row 60
$ws.Range("H60").Value = 26023.25
$ws.Range("J60").Value = 30000
$ws.Range("N60").Value = -31022
$ws.Range("L60").Value = 30000
# Row 123
$ws.Range("K123").Value = 65000
$ws.Range("M123").Value = -60100
$ws.Range("L123").Value = 100390
$ws.Range("I123").Value = 65000
$ws.Range("H123").Value = 88593.336
$ws.Range("J123").Value = 100390
$ws.Range("N123").Value = -110190
# Row 132
$ws.Range("I132").Value = 1949.25
$ws.Range("H132").Value = 2216.4783
$ws.Range("K132").Value = 5847.75
$ws.Range("M132").Value = -3317.75
# Row 134
$ws.Range("I134").Value = 1825.1578
$ws.Range("H134").Value = 2273.75
$ws.Range("M134").Value = -2940.4734
$ws.Range("K134").Value = 5475.4734

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("K4").Value = 5793096.6
$ws.Range("I4").Value = 1931032.2
$ws.Range("H4").Value = 10277351
$ws.Range("M4").Value = -5792984.6
# Row 11
$ws.Range("I11").Value = 500
$ws.Range("H11").Value = 500
$ws.Range("K11").Value = 1500
$ws.Range("M11").Value = -1360
# Row 51
$ws.Range("I51").Value = 2165
$ws.Range("L51").Value = 6510
$ws.Range("H51").Value = 2167.5
$ws.Range("J51").Value = 2170
$ws.Range("N51").Value = -7430
$ws.Range("K51").Value = 6495
$ws.Range("M51").Value = -6035
# Row 121
$ws.Range("H121").Value = 2112.4167
$ws.Range("J121").Value = 1671.1666
$ws.Range("N121").Value = -7633.4998
$ws.Range("L121").Value = 5013.4998
# Row 131
$ws.Range("H131").Value = 37553.195
$ws.Range("M131").Value = -830625.75
$ws.Range("K131").Value = 835665.75
$ws.Range("I131").Value = 278555.25

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("L70").Value = 10513.333
$ws.Range("I70").Value = 8577.223
$ws.Range("H70").Value = 9351.666999999999
$ws.Range("K70").Value = 8577.223
$ws.Range("J70").Value = 10513.333
$ws.Range("N70").Value = -11053.333
$ws.Range("M70").Value = -8307.223
# Row 73
$ws.Range("L73").Value = 10513.333
$ws.Range("I73").Value = 8577.223
$ws.Range("N73").Value = -12385.333
$ws.Range("M73").Value = -7641.223
$ws.Range("H73").Value = 9351.666999999999
$ws.Range("J73").Value = 10513.333
$ws.Range("K73").Value = 8577.223
# Row 107
$ws.Range("H107").Value = 424.73334
$ws.Range("J107").Value = 669
$ws.Range("N107").Value = -4509
$ws.Range("K107").Value = 145.57143
$ws.Range("M107").Value = 1774.42857
$ws.Range("L107").Value = 669
$ws.Range("I107").Value = 145.57143
# Row 126
$ws.Range("K126").Value = 7464.500100000001
$ws.Range("M126").Value = -4994.500100000001
$ws.Range("I126").Value = 2488.1667
$ws.Range("H126").Value = 3597.4348
# Row 132
$ws.Range("L132").Value = 9994.5
$ws.Range("H132").Value = 3903.2666
$ws.Range("J132").Value = 3331.5
$ws.Range("N132").Value = -15054.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 12
$ws.Range("H12").Value = 1500
$ws.Range("M12").Value = -1330
$ws.Range("K12").Value = 1500
$ws.Range("I12").Value = 1500
# Row 40
$ws.Range("L40").Value = 5594.4165
$ws.Range("I40").Value = 5052.7144
$ws.Range("H40").Value = 5394.8423
$ws.Range("J40").Value = 5594.4165
$ws.Range("N40").Value = -5866.4165
$ws.Range("K40").Value = 5052.7144
$ws.Range("M40").Value = -4916.7144
# Row 46
$ws.Range("H46").Value = 1972.6897
$ws.Range("J46").Value = 2351.1052
$ws.Range("N46").Value = -2727.1052
$ws.Range("M46").Value = -1065.7
$ws.Range("K46").Value = 1253.7
$ws.Range("L46").Value = 2351.1052
$ws.Range("I46").Value = 1253.7
# Row 82
$ws.Range("L82").Value = 884.5
$ws.Range("H82").Value = 1570.2
$ws.Range("J82").Value = 884.5
$ws.Range("N82").Value = -1606.5
# Row 85
$ws.Range("L85").Value = 884.5
$ws.Range("H85").Value = 1570.2
$ws.Range("J85").Value = 884.5
$ws.Range("N85").Value = -3380.5
# Row 136
$ws.Range("I136").Value = 4129.6665
$ws.Range("H136").Value = 3048.25
$ws.Range("M136").Value = -9838.999500000002
$ws.Range("K136").Value = 12388.9995

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Range("L54").Value = 15999
$ws.Range("I54").Value = 0
$ws.Range("H54").Value = 15999
$ws.Range("M54").ClearContents()
$ws.Range("J54").Value = 15999
$ws.Range("N54").Value = -17039
$ws.Range("K54").Value = 0
